$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.4379
$ws.Range("F3").Value = 0.009
$ws.Range("F4").Value = 0.4583
$ws.Range("B7").Value = 1.45909393024205
$ws.Range("C7").Value = 1.45909393024205
$ws.Range("D7").Value = 1.80254062412806
$ws.Range("E7").Value = 0.0120517545303688
$ws.Range("F7").Value = 0.1794
$ws.Range("B8").Value = 3.02226315113077
$ws.Range("C8").Value = 3.02226315113077
$ws.Range("D8").Value = 3.73365414919843
$ws.Range("E8").Value = 0.0249631451880308
$ws.Range("F8").Value = 0.0542
$ws.Range("B9").Value = 0.0246573092444615
$ws.Range("C9").Value = 0.0246573092444615
$ws.Range("D9").Value = 0.0304612339710419
$ws.Range("E9").Value = 0.000203663268165571
$ws.Range("F9").Value = 0.8681
$ws.Range("B10").Value = 116.562990671287
$ws.Range("C10").Value = 0.809465212995052
$ws.Range("E10").Value = 0.962781437013435
$ws.Range("B11").Value = 121.069005061905
$ws.Range("B12").Value = 0.599299757733328
$ws.Range("C12").Value = 0.599299757733328
$ws.Range("D12").Value = 0.951903536722973
$ws.Range("E12").Value = 0.00652121976240012
$ws.Range("F12").Value = 0.3337
$ws.Range("B13").Value = 0.166765352870391
$ws.Range("C13").Value = 0.166765352870391
$ws.Range("D13").Value = 0.26488335286599
$ws.Range("E13").Value = 0.00181464033780894
$ws.Range("F13").Value = 0.5972
$ws.Range("B14").Value = 0.474313210125965
$ws.Range("C14").Value = 0.474313210125965
$ws.Range("D14").Value = 0.753379951196707
$ws.Range("E14").Value = 0.00516119127286086
$ws.Range("F14").Value = 0.3838
$ws.Range("B15").Value = 90.659569782347
$ws.Range("C15").Value = 0.629580345710743
$ws.Range("E15").Value = 0.98650294862693
$ws.Range("B16").Value = 91.8999481030767
$ws.Range("F17").Value = 0.068
$ws.Range("F18").Value = 0.9969
$ws.Range("F19").Value = 0.574
$ws.Range("B22").Value = 0.385694596188176
$ws.Range("C22").Value = 0.385694596188176
$ws.Range("D22").Value = 1.90230852435001
$ws.Range("E22").Value = 0.0128631747690137
$ws.Range("F22").Value = 0.1721
$ws.Range("B23").Value = 0.0985394102351073
$ws.Range("C23").Value = 0.0985394102351073
$ws.Range("D23").Value = 0.486012409629954
$ws.Range("E23").Value = 0.00328635575405185
$ws.Range("F23").Value = 0.4826
$ws.Range("B24").Value = 0.304050651125929
$ws.Range("C24").Value = 0.304050651125929
$ws.Range("D24").Value = 1.49962729887155
$ws.Range("E24").Value = 0.0101402941672459
$ws.Range("F24").Value = 0.2304
$ws.Range("B25").Value = 29.1961167918724
$ws.Range("C25").Value = 0.202750811054669
$ws.Range("E25").Value = 0.973710175309689
$ws.Range("B26").Value = 29.9844014494216
$ws.Range("B27").Value = 0.00549445930836925
$ws.Range("C27").Value = 0.00549445930836925
$ws.Range("D27").Value = 0.0658861316173537
$ws.Range("E27").Value = 0.000436073834453364
$ws.Range("F27").Value = 0.8886
$ws.Range("B28").Value = 0.580507985055344
$ws.Range("C28").Value = 0.580507985055344
$ws.Range("D28").Value = 6.96108995657172
$ws.Range("E28").Value = 0.0460726577023305
$ws.Range("F28").Value = 0.0159
$ws.Range("B29").Value = 0.00520527723950737
$ws.Range("C29").Value = 0.00520527723950737
$ws.Range("D29").Value = 0.0624184404795942
$ws.Range("E29").Value = 0.000413122580008423
$ws.Range("F29").Value = 0.9015
$ws.Range("B30").Value = 12.0086294487621
$ws.Range("C30").Value = 0.0833932600608481
$ws.Range("E30").Value = 0.953078145883208
$ws.Range("B31").Value = 12.5998371703654
$ws.Range("F32").Value = 0.173
$ws.Range("F33").Value = 0.7062
$ws.Range("F34").Value = 0.9651
